# Cap nhat 1 so thu lat vat
# Slide 18, "Content Placeholder 2": re-split the existing bullet text into
# multiple runs (as PowerPoint's spell-checker does while editing), and add
# a new bullet paragraph at the end.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Paragraph 1: "Giao diện tab cơ bản" -----------------------------------
$para = $tr.Paragraphs(1, 1)
$para.Characters(1, 4).Text = "Giao"
$para.Characters(6, 4).Text = "diện"
$para.Characters(15, 2).Text = "cơ"
$para.Characters(18, 3).Text = "bản"

# --- Paragraph 2: "Hỗ trợ ngôn ngữ C, C++, C#, VB" -------------------------
$para = $tr.Paragraphs(2, 1)
$para.Characters(1, 2).Text = "Hỗ"
$para.Characters(4, 3).Text = "trợ"
$para.Characters(8, 4).Text = "ngôn"
$para.Characters(13, 3).Text = "ngữ"

# --- Paragraph 3: "Syntax highlight với từ khóa cố định" -------------------
$para = $tr.Paragraphs(3, 1)
$para.Characters(18, 3).Text = "với"
$para.Characters(22, 2).Text = "từ"
$para.Characters(25, 4).Text = "khóa"
$para.Characters(30, 2).Text = "cố"
$para.Characters(33, 4).Text = "định"

# --- Paragraph 4: "Auto complete với từ khóa cố định" ----------------------
$para = $tr.Paragraphs(4, 1)
$para.Characters(15, 3).Text = "với"
$para.Characters(19, 2).Text = "từ"
$para.Characters(22, 4).Text = "khóa"
$para.Characters(27, 2).Text = "cố"
$para.Characters(30, 4).Text = "định"

# --- Paragraph 5: "Code folding hoạt động với mức độ cơ bản" ---------------
$para = $tr.Paragraphs(5, 1)
$para.Characters(14, 4).Text = "hoạt"
$para.Characters(19, 4).Text = "động"
$para.Characters(24, 3).Text = "với"
$para.Characters(28, 3).Text = "mức"
$para.Characters(32, 2).Text = "độ"
$para.Characters(35, 2).Text = "cơ"
$para.Characters(38, 3).Text = "bản"

# --- Paragraph 6: "Document map hoạt động tốt" -----------------------------
$para = $tr.Paragraphs(6, 1)
$para.Characters(14, 4).Text = "hoạt"
$para.Characters(19, 4).Text = "động"
$para.Characters(24, 3).Text = "tốt"

# --- Paragraph 7: "Bookmark margin với Number margin hoạt động tốt" --------
$para = $tr.Paragraphs(7, 1)
$para.Characters(17, 3).Text = "với"
$para.Characters(35, 4).Text = "hoạt"
$para.Characters(40, 4).Text = "động"
$para.Characters(45, 3).Text = "tốt"

# --- New paragraph 8: "Auto Indenting và Brace Matching chạy hiệu quả" -----
$tr.InsertAfter("`rAuto Indenting va Brace Matching chay hieu qua")
$para = $tr.Paragraphs(8, 1)
$para.Characters(16, 2).Text = "và"
$para.Characters(34, 4).Text = "chạy"
$para.Characters(39, 4).Text = "hiệu"
$para.Characters(44, 3).Text = "quả"
